$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: the order in which brand-new text values are assigned below matches
# the order new entries were appended to the shared string table in the
# target workbook (dusky rockfish, northern rockfish, Sebastes_variabilis
# filename, the variabilis drive url, the polyspinis filename, year_start,
# year_end, the polyspinis drive url). Preserve this order so the shared
# string table indices line up.

# Row 5 common_name (new)
$ws.Range("C5").Value = "dusky rockfish"
# Row 6 common_name (new)
$ws.Range("C6").Value = "northern rockfish"
# Row 5 filename (new)
$ws.Range("D5").Value = "Data_Geostat_Sebastes_variabilis.rds"
# Row 6 url (new) -- this row's url is the "1hiaC..." drive link
$ws.Range("H6").Value = "https://drive.google.com/drive/folders/1hiaCSmB8vajir228x-Xo1zhfhdLi2ZgV"
# Row 6 filename (new)
$ws.Range("D6").Value = "Data_Geostat_Sebastes_polyspinis.RDS"
# New header columns for the year range
$ws.Range("F1").Value = "year_start"
$ws.Range("G1").Value = "year_end"
# Row 5 url (new) -- this row's url is the "1YRMh..." drive link
$ws.Range("H5").Value = "https://drive.google.com/drive/folders/1YRMhHOb9MMHa_YsKhxhd_1Sck7ehRY2Z"

# --- The existing "url" header and per-row url values shift from column F
#     to column H now that year_start/year_end occupy F and G. ---
$ws.Range("H1").Value = "url"
$ws.Range("H2").Value2 = $ws.Range("F2").Value2
$ws.Range("H3").Value2 = $ws.Range("F3").Value2
$ws.Range("H4").Value2 = $ws.Range("F4").Value2

# --- Fill in the year_start / year_end values for every data row ---
$ws.Range("F2").Value = 1990
$ws.Range("G2").Value = 2021
$ws.Range("F3").Value = 1990
$ws.Range("G3").Value = 2021
$ws.Range("F4").Value = 1990
$ws.Range("G4").Value = 2021
$ws.Range("F5").Value = 1990
$ws.Range("G5").Value = 2021
$ws.Range("F6").Value = 1984
$ws.Range("G6").Value = 2021

# --- Row 5 (Sebastes variabilis / dusky rockfish) remaining new cells ---
$ws.Range("A5").Value = 30152
$ws.Range("E5").Value = "RDS"

# --- Row 6 (Sebastes polyspinis / northern rockfish) remaining new cells ---
$ws.Range("A6").Value = 30420
$ws.Range("E6").Value = "RDS"

# --- Selection mirrors the saved view in the source workbook ---
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("F6").Select()
